# Add a new "2020" column (Q) to the indicator table, mirroring the
# formatting of the existing 2019 column (P), then fill in the figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the 2019 column (P4:P14) into the new column
# (Q4:Q14) - this is what a user does when extending the table with a new
# year: copy the previous year's column and paste its formats, then type
# the new numbers in.
$ws.Range("P4:P14").Copy()
$ws.Range("Q4:Q14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header year
$ws.Range("Q4").Value = 2020

# New data values for 2020
$ws.Range("Q5").Value = 4.4631700362051845
$ws.Range("Q6").Value = 22.107243650047039
$ws.Range("Q7").Value = 4.8469387755102042
$ws.Range("Q8").Value = 11.270912826533607
$ws.Range("Q9").Value = 8.2663605051664764
$ws.Range("Q10").Value = 9.0160381447767666
$ws.Range("Q11").Value = 2.7624309392265194
$ws.Range("Q12").Value = 1.1408815903197926
$ws.Range("Q13").Value = 1.7541111981205952
$ws.Range("Q14").Value = 3.6288232244686367

# Leave the selection where the author left it after typing in the new
# column - one row below the last data row.
$ws.Range("P15").Select()
